# Fruta / hortaliza, semanal
# Insert a new price-report row at row 21 (pushing existing rows 21-116 down
# to 22-117) and populate it with the new "Granada" (Wonderfull / Primera)
# observation for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 21:116 down by one, creating a blank row 21.
$ws.Rows("21:21").Insert()

# Fill in the new row with the reported values.
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "Vega Modelo de Temuco"
$ws.Range("C21").Value = "La Araucanía"
$ws.Range("D21").Value = 44677
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100104
$ws.Range("H21").Value = "Frutos de pepita"
$ws.Range("I21").Value = 100104001
$ws.Range("J21").Value = "Granada"
$ws.Range("K21").Value = "Wonderfull"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 65
$ws.Range("N21").Value = 12000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 12000
$ws.Range("Q21").Value = "$/bandeja 10 kilos"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 1200
$ws.Range("T21").Value = 10
